$wb = $excel.ActiveWorkbook

# --- Sheet "ENTREPRISES" (sheet1) ---
$ws1 = $wb.Worksheets.Item("ENTREPRISES")
$ws1.Range("P1").Value = "position"
$ws1.Range("N3").Value = 50.608265
$ws1.Range("O3").Value = 3.16046

# --- Sheet "LABOS" (sheet2) ---
$ws2 = $wb.Worksheets.Item("LABOS")
$ws2.Range("Q1").Value = "position"
$ws2.Range("O3").Value = 50.518746
$ws2.Range("P3").Value = 2.645622

# --- Sheet "FORMATIONS" (sheet3) ---
$ws3 = $wb.Worksheets.Item("FORMATIONS")
$ws3.Range("U1").Value = "position"
$ws3.Range("S4").Value = 50.934099
$ws3.Range("T4").Value = 1.808739
